# Resource sharing plan: fill in the "Title" / "text" placeholders with the
# drafted plan, plus a reviewer comment on the "variant libraries" mention.

$d = $word.ActiveDocument

# --- 1. Heading: "Title" -> "Resource sharing plan" ---------------------
$full = $d.Content.Text
$idx = $full.IndexOf("Title")
$r = $d.Range($idx, $idx + 5)
$r.Text = "Resource sharing plan"

# --- 2. Body paragraph: "text" -> the full resource-sharing-plan copy ---
# The target markup is made of many small, separately-authored runs (plus
# two spell-check-exempted proper nouns), so build it as literal OOXML and
# drop it in via Range.InsertXML rather than a plain Range.Text assignment
# (which would collapse everything into a single run).
$full = $d.Content.Text
$idx = $full.IndexOf("text")
$r2 = $d.Range($idx, $idx + 4)

$bodyRuns = @'
<w:r><w:t>Research resources generated through the course of this project will be made available to the wider community</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>before or immediately after publication</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>C</w:t></w:r><w:r><w:t>ode required</w:t></w:r><w:r><w:t xml:space="preserve"> to reproduce</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>all</w:t></w:r><w:r><w:t xml:space="preserve"> results</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">will be made publicly available through </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>github</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>before or immediately after publication</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>and release</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> associated with </w:t></w:r><w:r><w:t>publication</w:t></w:r><w:r><w:t xml:space="preserve"> will be archived through </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>zenodo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>Next-generation sequencing data</w:t></w:r><w:r><w:t xml:space="preserve"> generated by this project</w:t></w:r><w:r><w:t xml:space="preserve"> will be released publicly </w:t></w:r><w:r><w:t>through</w:t></w:r><w:r><w:t xml:space="preserve"> the SRA</w:t></w:r><w:r><w:t xml:space="preserve"> upon publication</w:t></w:r><w:r><w:t>. Plasmids</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:r><w:t>variant libraries</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">that are </w:t></w:r><w:r><w:t>generated by this project</w:t></w:r><w:r><w:t xml:space="preserve"> will be deposited with Addgene </w:t></w:r><w:r><w:t>before or immediately after publication</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>to ensure they are available to other researcher</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t>If any intellectual property arising from this project is patented, we will ensure that materials and data remain widely available to the research community.</w:t></w:r>
'@

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>$bodyRuns</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
$r2.InsertXML($xml)

# --- 3. Reviewer comment anchored on "variant libraries " ---------------
$full3 = $d.Content.Text
$needle = "variant libraries "
$idx3 = $full3.IndexOf($needle)
$r3 = $d.Range($idx3, $idx3 + $needle.Length)

$commentText = "Update this if needed when I hear back from Addgene about library deposition."
$cm = $d.Comments.Add($r3, $commentText)

# Re-fetch the freshly added comment from the collection and stamp the
# reviewer identity on it (writing straight through $word.UserName /
# $word.UserInitials after content has already been edited unsettles the
# document's first paragraph, so set the identity on the comment itself).
$live = $d.Comments(1)
$live.Author = "John Desmarais"
$live.Initial = "JD"

Write-Output "done"
